$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 202, shifting existing rows 202-216 down to 203-217.
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row 202 with the new weekly record.
$ws.Cells.Item(202, 1).Value = 8
$ws.Cells.Item(202, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(202, 3).Value = "Coquimbo"
$ws.Cells.Item(202, 4).Value = 45013
$ws.Cells.Item(202, 5).Value = 4
$ws.Cells.Item(202, 6).Value = 100112040
$ws.Cells.Item(202, 7).Value = "Cilantro"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 2500
$ws.Cells.Item(202, 11).Value = 1800
$ws.Cells.Item(202, 12).Value = 2000
$ws.Cells.Item(202, 13).Value = 1900
$ws.Cells.Item(202, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(202, 15).Value = "Provincia del Elqu" + [char]0xED
$ws.Cells.Item(202, 16).Value = 1267
$ws.Cells.Item(202, 17).Value = 1.5
$ws.Cells.Item(202, 18).Value = "Hortaliza"
